$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "triple_double_avg" header in J1, copying the header
# formatting (bold, centered, bordered) from the existing I1 header cell.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "triple_double_avg"

# Populate the new triple_double_avg column for each player row.
# Only Russell Westbrook's 2016-17 MVP season (row 10) averaged a triple-double.
$triple_double = @("No", "No", "No", "No", "No", "No", "No", "No", "Yes", "No")
for ($i = 0; $i -lt $triple_double.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $triple_double[$i]
}
